$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, E are never numeric-looking text, so a direct .Value assignment
# is safe. Column D frequently holds numeric-looking text (prices) that must stay
# TEXT (matches the source t="inlineStr" cells) instead of being coerced to a
# number by Excel type-inference. We force text entry by switching the cell to the
# "@" (Text) number format before the assignment, then ClearFormats() right after so
# no stray number-format/style is left behind on the cell (matches the target which
# carries no "s" style attribute on these cells).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.133.73"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.57%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.380.01"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.73%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.09"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.74%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.80"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.90%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.509"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.49%  "

# Row 8
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.91%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.21"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.48%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("D11").ClearFormats()

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.122"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.16%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.45"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.00%  "

# Row 14
$ws.Range("E14").Value = "  +1.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.750.36"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.378.35"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.37%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.807"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.68%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.136.21"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.70%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.17"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.03%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.33"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.10%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0888"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.15%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.59"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.63%  "

# Row 23
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.26"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.02%  "

# Row 24
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.15"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.18%  "

# Row 25
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.44"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.89%  "

# Row 26
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.89"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.79%  "

# Row 28
$ws.Range("E28").Value = "  +15.39%  "

# Row 29
$ws.Range("E29").Value = "  +1.23%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.49"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.20%  "

# Row 31
$ws.Range("E31").Value = "  -0.03%  "

# Row 32
$ws.Range("E32").Value = "  +2.56%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0739"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.42%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.14"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.96%  "

# Row 35
$ws.Range("E35").Value = "  +6.83%  "

# Row 36
$ws.Range("E36").Value = "  +3.33%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.29"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.40%  "

# Row 38
$ws.Range("E38").Value = "  -1.21%  "

# Row 39
$ws.Range("E39").Value = "  +4.66%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.56"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +13.88%  "

# Row 41
$ws.Range("E41").Value = "  +0.31%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.955.90"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.43%  "

# Row 43
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.12"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -37.20%  "

# Row 44
$ws.Range("E44").Value = "  +0.82%  "

# Row 45
$ws.Range("E45").Value = "  +2.22%  "

# Row 46
$ws.Range("E46").Value = "  +0.70%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.21"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -10.79%  "

# Row 48
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.607.00"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.48%  "

# Row 49
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.75"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.29%  "

# Row 50
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.51"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.94%  "

# Row 51
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.76"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.90%  "
